$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (report volume/number and week-covering date range)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  15"
$ws.Range("C9").Value = "Report Covering the Week  4/10/2023  Through  4/16/2023"

# ---------------------------------------------------------------------------
# Cells that change from a numeric value to the text placeholder "0"
# (used in this report to denote a blank/zero count cell). NumberFormat is
# temporarily forced to Text ("@") so Excel stores the value as a literal
# string instead of re-parsing "0" back into a number, then restored to
# General to match the look of the other placeholder cells in the sheet.
# ---------------------------------------------------------------------------
$numToTextCells = @("C16", "C17", "C27", "F28", "F29")
foreach ($addr in $numToTextCells) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = "0"
    $c.NumberFormat = "General"
}

# ---------------------------------------------------------------------------
# Cells that change from the text placeholder to an actual numeric value.
# Copy the number format from a sibling cell that already carries the
# desired numeric display format before assigning the new numeric value.
# ---------------------------------------------------------------------------
$ws.Range("D27").NumberFormat = $ws.Range("D18").NumberFormat
$ws.Range("D27").Value = 1

$ws.Range("E27").NumberFormat = $ws.Range("E18").NumberFormat
$ws.Range("E27").Value = -100

# ---------------------------------------------------------------------------
# Remaining plain numeric value updates (style/number format unchanged)
# ---------------------------------------------------------------------------
    $ws.Range("F15").Value = 3
    $ws.Range("I15").Value = 3
    $ws.Range("M15").Value = 50
    $ws.Range("N15").Value = 50
    $ws.Range("D16").Value = 5
    $ws.Range("E16").Value = -100
    $ws.Range("F16").Value = 4
    $ws.Range("G16").Value = 9
    $ws.Range("H16").Value = -55.555555555555
    $ws.Range("J16").Value = 29
    $ws.Range("K16").Value = -31.03448275862
    $ws.Range("L16").Value = 400
    $ws.Range("N16").Value = -81.981981981982
    $ws.Range("E17").Value = -100
    $ws.Range("F17").Value = 10
    $ws.Range("H17").Value = 150
    $ws.Range("J17").Value = 21
    $ws.Range("K17").Value = 23.809523809523
    $ws.Range("L17").Value = 100
    $ws.Range("M17").Value = 44.444444444444
    $ws.Range("N17").Value = -23.529411764705
    $ws.Range("C18").Value = 5
    $ws.Range("D18").Value = 4
    $ws.Range("E18").Value = 25
    $ws.Range("F18").Value = 20
    $ws.Range("G18").Value = 21
    $ws.Range("H18").Value = -4.761904761904
    $ws.Range("I18").Value = 102
    $ws.Range("J18").Value = 96
    $ws.Range("K18").Value = 6.25
    $ws.Range("L18").Value = 50
    $ws.Range("M18").Value = 29.113924050632
    $ws.Range("N18").Value = -64.583333333333
    $ws.Range("C19").Value = 15
    $ws.Range("D19").Value = 10
    $ws.Range("E19").Value = 50
    $ws.Range("F19").Value = 62
    $ws.Range("G19").Value = 46
    $ws.Range("H19").Value = 34.782608695652
    $ws.Range("I19").Value = 222
    $ws.Range("J19").Value = 158
    $ws.Range("K19").Value = 40.506329113924
    $ws.Range("L19").Value = 143.956043956044
    $ws.Range("M19").Value = 91.379310344827
    $ws.Range("N19").Value = 52.054794520547
    $ws.Range("D20").Value = 1
    $ws.Range("E20").Value = 0
    $ws.Range("G20").Value = 6
    $ws.Range("H20").Value = 16.666666666666
    $ws.Range("I20").Value = 41
    $ws.Range("J20").Value = 31
    $ws.Range("K20").Value = 32.258064516129
    $ws.Range("L20").Value = 127.777777777778
    $ws.Range("M20").Value = 2.5
    $ws.Range("N20").Value = -95.543478260869
    $ws.Range("C21").Value = 22
    $ws.Range("D21").Value = 21
    $ws.Range("E21").Value = 4.761904761904
    $ws.Range("F21").Value = 106
    $ws.Range("G21").Value = 86
    $ws.Range("H21").Value = 23.255813953488
    $ws.Range("I21").Value = 415
    $ws.Range("J21").Value = 335
    $ws.Range("K21").Value = 23.880597014925
    $ws.Range("L21").Value = 113.917525773196
    $ws.Range("M21").Value = 49.280575539568
    $ws.Range("N21").Value = -72.37017310253
    $ws.Range("C24").Value = 13
    $ws.Range("D24").Value = 18
    $ws.Range("E24").Value = -27.777777777777
    $ws.Range("F24").Value = 37
    $ws.Range("G24").Value = 60
    $ws.Range("H24").Value = -38.333333333333
    $ws.Range("I24").Value = 155
    $ws.Range("J24").Value = 235
    $ws.Range("K24").Value = -34.042553191489
    $ws.Range("L24").Value = 6.164383561643
    $ws.Range("M24").Value = 9.154929577464
    $ws.Range("D25").Value = 3
    $ws.Range("E25").Value = -33.333333333333
    $ws.Range("F25").Value = 14
    $ws.Range("H25").Value = 0
    $ws.Range("I25").Value = 59
    $ws.Range("J25").Value = 63
    $ws.Range("K25").Value = -6.349206349206
    $ws.Range("L25").Value = 73.529411764705
    $ws.Range("M25").Value = 51.282051282051
    $ws.Range("F26").Value = 4
    $ws.Range("I26").Value = 4
    $ws.Range("L26").Value = 300
    $ws.Range("J27").Value = 8
    $ws.Range("K27").Value = -37.5
    $ws.Range("L27").Value = 25
    $ws.Range("G28").Value = 1
    $ws.Range("H28").Value = -100
    $ws.Range("G29").Value = 1
    $ws.Range("H29").Value = -100
